$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete row 45 first (the stray duplicate row), then row 32
# (FC_NON_INT_EXP_GROSS_REV_BNK / its formula), so row indices for the
# not-yet-deleted row stay valid.
$ws.Rows.Item(45).Delete()
$ws.Rows.Item(32).Delete()

$ws.Range("A32:XFD32").Select()
